# Apply the target changes to the presentation:
#  1) Retarget the table style used by the three data tables (slides 14-16)
#     from the custom "Table_0" style {E76F2C34-F49C-4B5A-AC4B-6A491A10CF90}
#     to {82B68193-6186-4F62-B64E-E4511BABF635}.
#  2) Swap the deck's colour theme from the "Integral"/"Red Violet" palette
#     back to the standard "Office" palette.

$p = $ppt.ActivePresentation

$newTableStyleId = "{82B68193-6186-4F62-B64E-E4511BABF635}"
$targetSlides = @(14, 15, 16)

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# Restore the standard Office colour palette on the presentation's theme.
$officeColors = @(
    0,          # dk1 - 000000
    16777215,   # lt1 - FFFFFF
    6968388,    # dk2 - 44546A
    15132391,   # lt2 - E7E6E6
    13998939,   # accent1 - 5B9BD5
    3243501,    # accent2 - ED7D31
    10855845,   # accent3 - A5A5A5
    49407,      # accent4 - FFC000
    12874308,   # accent5 - 4472C4
    4697456,    # accent6 - 70AD47
    12673797,   # hlink - 0563C1
    7491477     # folHlink - 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
